$d = $word.ActiveDocument
$d.Content.Find.Execute("by simple integration of the", $true, $false, $false, $false, $false, $true, 1, $false, "by step integration of the", 2)
